$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 173, pushing existing rows 173:213 down to 174:214
$ws.Rows("173:173").Insert()

# Populate the newly inserted row 173 with the new data record
$ws.Range("A173").Value = 10
$ws.Range("B173").Value = "Vega Modelo de Temuco"
$ws.Range("C173").Value = "La Araucanía"
$ws.Range("D173").Value = 44642
$ws.Range("E173").Value = 9
$ws.Range("F173").Value = "Fruta"
$ws.Range("G173").Value = 100102
$ws.Range("H173").Value = "Cítricos"
$ws.Range("I173").Value = 100102006
$ws.Range("J173").Value = "Pomelo"
$ws.Range("K173").Value = "Start Ruby"
$ws.Range("L173").Value = "Primera"
$ws.Range("M173").Value = 80
$ws.Range("N173").Value = 14000
$ws.Range("O173").Value = 14000
$ws.Range("P173").Value = 14000
$ws.Range("Q173").Value = "$/bandeja 15 kilos granel"
$ws.Range("R173").Value = "Región de O'Higgins"
$ws.Range("S173").Value = 933
$ws.Range("T173").Value = 15
